$d = $word.ActiveDocument

# 1) Append four soft hyphens to the end of the "It should not cover..." paragraph.
$softHyphen = [char]31
$targetText = "It should not cover the whole screen.(it should cover the space smaller than other section like it look in the picture)"
$found = $d.Content.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, ($targetText + $softHyphen + $softHyphen + $softHyphen + $softHyphen), 2)

# 2) Remove the "img_03, img_04, img_05 and img_06." text from the "Include content with its
#    image   img_03, ..." paragraph, leaving the trailing spaces behind.
$found = $d.Content.Find.Execute("   img_03, img_04, img_05 and img_06.", $true, $false, $false, $false, $false, $true, 1, $false, "   ", 2)

# 3) Re-append the removed text to the end of the "...hover" paragraph, and insert a new
#    empty paragraph right after it (so Section_08 stays on its own paragraph).
$found = $d.Content.Find.Execute("hover", $true, $false, $false, $false, $false, $true, 1, $false, "hover img_03, img_04, img_05 and img_06.^p", 2)
